$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff/Handback Datetime strings in row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 21:19:57"
$wsZhCn.Range("H2").Value = "2016-03-23 21:20:39"

# de-de sheet: update Correspond Handoff/Handback Datetime strings in row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 21:20:01"
$wsDeDe.Range("H2").Value = "2016-03-23 21:20:46"
